$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 119 ---
# A119: date serial changes from 45460.6493865741 to 45460.2916666667
$ws.Range("A119").Value2 = 45460.2916666667
# D119: low value changes from 1.57000005245209 to 1.55499994754791
$ws.Range("D119").Value2 = 1.55499994754791

# --- Append new row 120 (same layout/style as row 119) ---
# Copy A119's format (date/time number format) onto A120 before writing its value
$ws.Range("A119").Copy()
$ws.Range("A120").PasteSpecial(-4122)

$ws.Range("A120").Value2 = 45461.6457175926
$ws.Range("B120").Value2 = 6000
$ws.Range("C120").Value2 = 1.58000004291534
$ws.Range("D120").Value2 = 1.56500005722046
$ws.Range("E120").Value2 = 1.57000005245209
$ws.Range("F120").Value2 = 1.58000004291534

# G120 and H120 must be stored as text (shared strings), not auto-detected numbers
$ws.Range("G120").NumberFormat = "@"
$ws.Range("G120").Value = "1.58000004291534"
$ws.Range("G120").Style = "Normal"

$ws.Range("H120").NumberFormat = "@"
$ws.Range("H120").Value = "SMN.MI"
$ws.Range("H120").Style = "Normal"
